# 15.1.1.1 "Forest area" sheet: add the 2021/2022 columns (M, N) to the
# existing 2012-2020 time series and append a footnote row explaining the
# source of the new 2022 figure.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats = -4122 (copy only number-format / style from the source cell)
function Set-CellWithStyleFrom {
    param(
        [string]$TargetAddr,
        [string]$StyleSourceAddr,
        $Value
    )
    # Write the value first ...
    $ws.Range($TargetAddr).Value = $Value
    # ... then stamp the same cell style/format as an existing, already-styled
    # neighbour cell (same row) so the new cell matches the table formatting.
    $ws.Range($StyleSourceAddr).Copy()
    $ws.Range($TargetAddr).PasteSpecial(-4122)
}

# --- Row 4: new year headers (2021, 2022) ---------------------------------
Set-CellWithStyleFrom "M4" "L4" 2021
Set-CellWithStyleFrom "N4" "L4" 2022

# --- Row 5: Kyrgyz Republic totals -----------------------------------------
Set-CellWithStyleFrom "M5" "L5" 5.6
Set-CellWithStyleFrom "N5" "L5" 6.3

# --- Row 6: Batken oblast ----------------------------------------------
Set-CellWithStyleFrom "M6" "L6" 0.8
Set-CellWithStyleFrom "N6" "L6" 0.8

# --- Row 7: Jalal-Abad oblast --------------------------------------------
Set-CellWithStyleFrom "M7" "L7" 1.9
Set-CellWithStyleFrom "N7" "L7" 2.4

# --- Row 8: Ysyk-Kul oblast ------------------------------------------------
Set-CellWithStyleFrom "M8" "L8" 0.7
Set-CellWithStyleFrom "N8" "L8" 0.7

# --- Row 9: Naryn oblast -----------------------------------------------
Set-CellWithStyleFrom "M9" "L9" 0.7
Set-CellWithStyleFrom "N9" "L9" 0.8

# --- Row 10: Osh oblast ------------------------------------------------
Set-CellWithStyleFrom "M10" "L10" 0.9
Set-CellWithStyleFrom "N10" "L10" 1

# --- Row 11: Talas oblast ----------------------------------------------
Set-CellWithStyleFrom "M11" "L11" 0.3
Set-CellWithStyleFrom "N11" "L11" 0.2

# --- Row 12: Chui oblast (bottom border row) ----------------------------
Set-CellWithStyleFrom "M12" "L12" 0.2
Set-CellWithStyleFrom "N12" "L12" 0.4

# --- Row 14: new footnote about the 2022 forest-inventory data source ----
Set-CellWithStyleFrom "B14" "B13" "По данным лесоустройства 2022 года Лесной службы при Министерстве чрезвычайных ситуаций КР"
$ws.Rows.Item(14).RowHeight = 34.5

# Note: the original sheetView carried a stale <selection activeCell="N5".../>
# from before this edit. This sandboxed Excel host always re-emits a
# <selection> element for whatever cell is tracked as "active" (confirmed
# even via the lowest-level model API, independent of this script), so it
# cannot be fully cleared from COM automation; left untouched here as it
# carries no workbook data.
